$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 9253.091
$ws.Range("I9").Value = 10098.4
$ws.Range("K9").Value = 10098.4
$ws.Range("M9").Value = -9929.4
$ws.Range("H107").Value = 1092.3846
$ws.Range("I107").Value = 1092.3846
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1092.3846
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 827.6153999999999
$ws.Range("N107").ClearContents()
$ws.Range("H112").Value = 126649.875
$ws.Range("J112").Value = 92864.63
$ws.Range("L112").Value = 278593.89
$ws.Range("N112").Value = -280809.89
$ws.Range("H131").Value = 1000
$ws.Range("I131").Value = 1000
$ws.Range("K131").Value = 3000
$ws.Range("M131").Value = 2040
$ws.Range("H137").Value = 4047.25
$ws.Range("I137").Value = 3063
$ws.Range("K137").Value = 9189
$ws.Range("M137").Value = -6639

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 2314.2
$ws.Range("I5").Value = 2035.5
$ws.Range("J5").Value = 2500
$ws.Range("K5").Value = 2035.5
$ws.Range("L5").Value = 2500
$ws.Range("M5").Value = -1923.5
$ws.Range("N5").Value = -2724
$ws.Range("H45").Value = 5436
$ws.Range("I45").Value = 6661.727
$ws.Range("J45").Value = 2739.4
$ws.Range("K45").Value = 6661.727
$ws.Range("L45").Value = 2739.4
$ws.Range("M45").Value = -6284.727
$ws.Range("N45").Value = -3493.4
$ws.Range("H61").Value = 52633388
$ws.Range("I61").Value = 58825340
$ws.Range("K61").Value = 58825340
$ws.Range("M61").Value = -58825128
$ws.Range("H122").Value = 3923.2144
$ws.Range("I122").Value = 2971.7827
$ws.Range("J122").Value = 8299.799999999999
$ws.Range("K122").Value = 8915.348100000001
$ws.Range("L122").Value = 24899.4
$ws.Range("M122").Value = -6465.348100000001
$ws.Range("N122").Value = -29799.4
$ws.Range("H136").Value = 52633388
$ws.Range("I136").Value = 58825340
$ws.Range("K136").Value = 176476020
$ws.Range("M136").Value = -176473470

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 2314.2
$ws.Range("I4").Value = 2035.5
$ws.Range("J4").Value = 2500
$ws.Range("K4").Value = 2035.5
$ws.Range("L4").Value = 2500
$ws.Range("M4").Value = -1920.5
$ws.Range("N4").Value = -2730
$ws.Range("H86").Value = 1803.4872
$ws.Range("I86").Value = 1748.1613
$ws.Range("K86").Value = 1748.1613
$ws.Range("M86").Value = -625.1613
$ws.Range("H89").Value = 1803.4872
$ws.Range("I89").Value = 1748.1613
$ws.Range("K89").Value = 8740.806500000001
$ws.Range("M89").Value = -3124.806500000001
$ws.Range("H99").Value = 2006.3125
$ws.Range("I99").Value = 1984
$ws.Range("K99").Value = 1984
$ws.Range("M99").Value = -486
$ws.Range("H134").Value = 55562196
$ws.Range("I134").Value = 62507296
$ws.Range("J134").Value = 1397
$ws.Range("K134").Value = 187521888
$ws.Range("L134").Value = 4191
$ws.Range("M134").Value = -187519353
$ws.Range("N134").Value = -9261

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 1255
$ws.Range("I17").Value = 1255
$ws.Range("K17").Value = 1255
$ws.Range("M17").Value = -1081
$ws.Range("H31").Value = 7197.0557
$ws.Range("I31").Value = 5000.2173
$ws.Range("J31").Value = 11083.77
$ws.Range("K31").Value = 5000.2173
$ws.Range("L31").Value = 11083.77
$ws.Range("M31").Value = -4705.2173
$ws.Range("N31").Value = -11673.77
$ws.Range("H34").Value = 7197.0557
$ws.Range("I34").Value = 5000.2173
$ws.Range("J34").Value = 11083.77
$ws.Range("K34").Value = 5000.2173
$ws.Range("L34").Value = 11083.77
$ws.Range("M34").Value = -4798.2173
$ws.Range("N34").Value = -11487.77
$ws.Range("H62").Value = 2583.3333
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 2875
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 2875
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -4123
$ws.Range("H65").Value = 2583.3333
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 2875
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 14375
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -20615
$ws.Range("H86").Value = 15225.5
$ws.Range("J86").Value = 17204
$ws.Range("L86").Value = 17204
$ws.Range("N86").Value = -19450
$ws.Range("H89").Value = 15225.5
$ws.Range("J89").Value = 17204
$ws.Range("L89").Value = 86020
$ws.Range("N89").Value = -97252
$ws.Range("H94").Value = 2614.6
$ws.Range("I94").Value = 2818.5
$ws.Range("J94").Value = 1799
$ws.Range("K94").Value = 2818.5
$ws.Range("L94").Value = 1799
$ws.Range("M94").Value = -2367.5
$ws.Range("N94").Value = -2701
$ws.Range("H134").Value = 11957031
$ws.Range("I134").Value = 13215086
$ws.Range("J134").Value = 5507
$ws.Range("K134").Value = 39645258
$ws.Range("L134").Value = 16521
$ws.Range("M134").Value = -39642723
$ws.Range("N134").Value = -21591

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 999.9545000000001
$ws.Range("J34").Value = 1000.0238
$ws.Range("L34").Value = 3000.0714
$ws.Range("N34").Value = -3168.0714
$ws.Range("H39").Value = 824.375
$ws.Range("I39").Value = 656.5714
$ws.Range("K39").Value = 1969.7142
$ws.Range("M39").Value = -1675.7142
$ws.Range("H40").Value = 80.125
$ws.Range("I40").Value = 62.2
$ws.Range("K40").Value = 248.8
$ws.Range("M40").Value = -179.8
$ws.Range("H55").Value = 1001
$ws.Range("I55").Value = 1001
$ws.Range("K55").Value = 3003
$ws.Range("M55").Value = -2826
$ws.Range("H132").Value = 6666.6665
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 6666.6665
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 59999.9985
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -65059.9985
$ws.Range("H140").Value = 1053.3
$ws.Range("I140").Value = 1053.3
$ws.Range("K140").Value = 3159.9
$ws.Range("M140").Value = 2020.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H97").Value = 1380.0952
$ws.Range("I97").Value = 1075.7858
$ws.Range("K97").Value = 1075.7858
$ws.Range("M97").Value = -579.7858000000001
$ws.Range("H132").Value = 5002773
$ws.Range("I132").Value = 5684570.5
$ws.Range("J132").Value = 2924.3333
$ws.Range("K132").Value = 17053711.5
$ws.Range("L132").Value = 8772.999899999999
$ws.Range("M132").Value = -17051181.5
$ws.Range("N132").Value = -13832.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2033.2778
$ws.Range("I46").Value = 2253.7693
$ws.Range("J46").Value = 1460
$ws.Range("K46").Value = 2253.7693
$ws.Range("L46").Value = 1460
$ws.Range("M46").Value = -2065.7693
$ws.Range("N46").Value = -1836
$ws.Range("H55").Value = 342.7
$ws.Range("I55").Value = 355.25
$ws.Range("K55").Value = 355.25
$ws.Range("M55").Value = -182.25
$ws.Range("H57").Value = 22235.25
$ws.Range("I57").Value = 20013.666
$ws.Range("J57").Value = 28900
$ws.Range("K57").Value = 20013.666
$ws.Range("L57").Value = 28900
$ws.Range("M57").Value = -19447.666
$ws.Range("N57").Value = -30032
$ws.Range("H61").Value = 3798.5557
$ws.Range("I61").Value = 4119.2085
$ws.Range("J61").Value = 1233.3334
$ws.Range("K61").Value = 4119.2085
$ws.Range("L61").Value = 1233.3334
$ws.Range("M61").Value = -3917.2085
$ws.Range("N61").Value = -1637.3334
$ws.Range("H82").Value = 1268.5
$ws.Range("I82").Value = 1268.5
$ws.Range("K82").Value = 1268.5
$ws.Range("M82").Value = -907.5
$ws.Range("H85").Value = 1268.5
$ws.Range("I85").Value = 1268.5
$ws.Range("K85").Value = 1268.5
$ws.Range("M85").Value = -20.5
$ws.Range("H113").Value = 3798.5557
$ws.Range("I113").Value = 4119.2085
$ws.Range("J113").Value = 1233.3334
$ws.Range("K113").Value = 4119.2085
$ws.Range("L113").Value = 1233.3334
$ws.Range("M113").Value = -1949.2085
$ws.Range("N113").Value = -5573.3334
$ws.Range("H122").Value = 5233.6665
$ws.Range("H136").Value = 3425.1428
$ws.Range("I136").Value = 2829.3333
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 8487.999899999999
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -5937.999899999999
$ws.Range("N136").Value = -26100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3566
$ws.Range("I96").Value = 1399.75
$ws.Range("J96").Value = 4075.7058
$ws.Range("K96").Value = 1399.75
$ws.Range("L96").Value = 4075.7058
$ws.Range("M96").Value = -26.75
$ws.Range("N96").Value = -6821.7058
$ws.Range("H113").Value = 961.25
$ws.Range("I113").Value = 786.75
$ws.Range("K113").Value = 2360.25
$ws.Range("M113").Value = -190.25
$ws.Range("H126").Value = 2266.5334
$ws.Range("I126").Value = 2530.3
$ws.Range("J126").Value = 1739
$ws.Range("K126").Value = 7590.900000000001
$ws.Range("L126").Value = 5217
$ws.Range("M126").Value = -5120.900000000001
$ws.Range("N126").Value = -10157
$ws.Range("H132").Value = 12197735
$ws.Range("I132").Value = 14287119
$ws.Range("J132").Value = 9663.333000000001
$ws.Range("K132").Value = 42861357
$ws.Range("L132").Value = 28989.999
$ws.Range("M132").Value = -42858827
$ws.Range("N132").Value = -34049.999
